{"js": "// Update the worksheet date and every division-problem cell to the new\n// values (output regenerated at commit c986bee).\nconst replacements = [\n  [\"2024-09-03 Tuesday\", \"2024-09-04 Wednesday\"],\n  [\"113\u00f77=16, 1\", \"412\u00f78=51, 4\"],\n  [\"500\u00f74=125, 0\", \"120\u00f76=20, 0\"],\n  [\"678\u00f73=226, 0\", \"145\u00f72=72, 1\"],\n  [\"145\u00f76=24, 1\", \"262\u00f72=131, 0\"],\n  [\"400\u00f72=200, 0\", \"646\u00f74=161, 2\"],\n  [\"550\u00f74=137, 2\", \"615\u00f74=153, 3\"],\n  [\"510\u00f75=102, 0\", \"342\u00f74=85, 2\"],\n  [\"816\u00f72=408, 0\", \"906\u00f76=151, 0\"],\n  [\"664\u00f74=166, 0\", \"791\u00f78=98, 7\"],\n  [\"334\u00f79=37, 1\", \"310\u00f77=44, 2\"],\n  [\"706\u00f75=141, 1\", \"883\u00f76=147, 1\"],\n  [\"383\u00f79=42, 5\", \"449\u00f78=56, 1\"],\n  [\"856\u00f74=214, 0\", \"366\u00f75=73, 1\"],\n  [\"176\u00f76=29, 2\", \"409\u00f79=45, 4\"],\n  [\"778\u00f79=86, 4\", \"273\u00f78=34, 1\"],\n  [\"370\u00f79=41, 1\", \"771\u00f72=385, 1\"],\n  [\"119\u00f78=14, 7\", \"242\u00f78=30, 2\"],\n  [\"355\u00f78=44, 3\", \"665\u00f76=110, 5\"],\n  [\"158\u00f78=19, 6\", \"139\u00f79=15, 4\"],\n  [\"962\u00f72=481, 0\", \"752\u00f78=94, 0\"],\n  [\"761\u00f77=108, 5\", \"407\u00f78=50, 7\"],\n  [\"519\u00f73=173, 0\", \"857\u00f74=214, 1\"],\n  [\"595\u00f73=198, 1\", \"661\u00f76=110, 1\"],\n  [\"350\u00f72=175, 0\", \"726\u00f72=363, 0\"],\n  [\"490\u00f78=61, 2\", \"260\u00f78=32, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every division-problem cell to the new\n# values (output regenerated at commit c986bee).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-09-03 Tuesday', '2024-09-04 Wednesday'),\n    @('113\u00f77=16, 1', '412\u00f78=51, 4'),\n    @('500\u00f74=125, 0', '120\u00f76=20, 0'),\n    @('678\u00f73=226, 0', '145\u00f72=72, 1'),\n    @('145\u00f76=24, 1', '262\u00f72=131, 0'),\n    @('400\u00f72=200, 0', '646\u00f74=161, 2'),\n    @('550\u00f74=137, 2', '615\u00f74=153, 3'),\n    @('510\u00f75=102, 0', '342\u00f74=85, 2'),\n    @('816\u00f72=408, 0', '906\u00f76=151, 0'),\n    @('664\u00f74=166, 0', '791\u00f78=98, 7'),\n    @('334\u00f79=37, 1', '310\u00f77=44, 2'),\n    @('706\u00f75=141, 1', '883\u00f76=147, 1'),\n    @('383\u00f79=42, 5', '449\u00f78=56, 1'),\n    @('856\u00f74=214, 0', '366\u00f75=73, 1'),\n    @('176\u00f76=29, 2', '409\u00f79=45, 4'),\n    @('778\u00f79=86, 4', '273\u00f78=34, 1'),\n    @('370\u00f79=41, 1', '771\u00f72=385, 1'),\n    @('119\u00f78=14, 7', '242\u00f78=30, 2'),\n    @('355\u00f78=44, 3', '665\u00f76=110, 5'),\n    @('158\u00f78=19, 6', '139\u00f79=15, 4'),\n    @('962\u00f72=481, 0', '752\u00f78=94, 0'),\n    @('761\u00f77=108, 5', '407\u00f78=50, 7'),\n    @('519\u00f73=173, 0', '857\u00f74=214, 1'),\n    @('595\u00f73=198, 1', '661\u00f76=110, 1'),\n    @('350\u00f72=175, 0', '726\u00f72=363, 0'),\n    @('490\u00f78=61, 2', '260\u00f78=32, 4'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
